$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.165667057037354
$ws.Range("B1").Value = 1.335435748100281
$ws.Range("C1").Value = 1.669107675552368
$ws.Range("D1").Value = 3.268141269683838
$ws.Range("E1").Value = 15
